$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, matching the style of the other
# header cells (bold / bordered / centered) by copying E1's format over.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Per-row timestamps recorded when each panel entry was processed.
$timestamps = @(
    "2021-10-05 13:40:37.023387",
    "2021-10-05 13:40:37.023398",
    "2021-10-05 13:40:37.023402",
    "2021-10-05 13:40:37.023406",
    "2021-10-05 13:40:37.023409",
    "2021-10-05 13:40:37.023412",
    "2021-10-05 13:40:37.023415",
    "2021-10-05 13:40:37.023418",
    "2021-10-05 13:40:37.023422",
    "2021-10-05 13:40:37.023425",
    "2021-10-05 13:40:37.023427",
    "2021-10-05 13:40:37.023430",
    "2021-10-05 13:40:37.023433",
    "2021-10-05 13:40:37.023436",
    "2021-10-05 13:40:37.023439",
    "2021-10-05 13:40:37.023442",
    "2021-10-05 13:40:37.023445"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

Write-Output "done"
